$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the data held in row 2 ("2019年B") and row 3 ("2019年C") ---
# Capture current row 2 and row 3 values (columns A:P) before overwriting.
$row2Values = $ws.Range("A2:P2").Value()
$row3Values = $ws.Range("A3:P3").Value()

$ws.Range("A2:P2").Value = $row3Values
$ws.Range("A3:P3").Value = $row2Values

# --- Step 2: remove the duplicated (non-cumulative) columns Q:AE entirely ---
$ws.Range("Q1:AE1").EntireColumn.Delete()
